# Heterogeneity table (t5_robust_diff): the outcome rows are reshuffled
# so the "Contraceptive knowledge" block moves up to sit right after the
# "Male condom" block (and swaps with "Female condom" block), and the
# "Male condom attitudes index" / "Used male condom..." rows swap order.
# Every data row (A:I) keeps its own RD / 95% CI values - only row
# position changes. Easiest robust way to express that via COM is to
# rewrite A3:I13 with the final, reordered content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final row order (outcome label + RD/95% CI for each of the 4 model
# columns: OLS Clustered SEs, OLS Pooled, GEE, HLM Ward Intercepts).
$rows = @(
    @('Male condom attitudes index', '0.011', '(-0.103, 0.125)', '0.014', '(-0.101, 0.130)', '0.011', '(-0.097, 0.119)', '0.011', '(-0.099, 0.120)'),
    @('Used male condom at most recent sex (0/1)', '-0.002', '(-0.046, 0.041)', '-0.006', '(-0.047, 0.036)', '-0.002', '(-0.044, 0.040)', '-0.002', '(-0.044, 0.039)'),
    @('Contraceptive knowledge index', '0.029', '(-0.068, 0.125)', '0.059', '(-0.061, 0.179)', '0.029', '(-0.071, 0.129)', '0.029', '(-0.072, 0.129)'),
    @('Modern contraceptive methods known (n)', '-0.020', '(-0.327, 0.287)', '-0.040', '(-0.386, 0.305)', '-0.015', '(-0.306, 0.276)', '-0.015', '(-0.307, 0.277)'),
    @('Discussed contraceptive use with recent partner (0/1)', '-0.012', '(-0.065, 0.040)', '-0.021', '(-0.078, 0.035)', '-0.011', '(-0.064, 0.042)', '-0.011', '(-0.064, 0.042)'),
    @('Can identify a female condom (0/1)', '0.048+', '(0.001, 0.096)', '0.057+', '(0.001, 0.112)', '0.048+', '(-0.003, 0.099)', '0.048+', '(-0.003, 0.099)'),
    @('Would be willing to try a female condom (0/1)', '-0.006', '(-0.054, 0.043)', '0.005', '(-0.056, 0.066)', '-0.006', '(-0.051, 0.040)', '-0.006', '(-0.051, 0.040)'),
    @('Female condom attitudes index', '-0.000', '(-0.111, 0.111)', '0.068', '(-0.063, 0.198)', '0.002', '(-0.102, 0.107)', '0.003', '(-0.104, 0.109)'),
    @('Has ever used a female condom (0/1)', '0.013', '(-0.005, 0.032)', '0.019', '(-0.005, 0.042)', '0.013', '(-0.004, 0.031)', '0.013', '(-0.007, 0.034)'),
    @('Used a female condom in last 6 months (0/1)', '0.007', '(-0.004, 0.018)', '0.011', '(-0.003, 0.024)', '0.007', '(-0.005, 0.019)', '0.007', '(-0.006, 0.020)'),
    @('Used a female condom at most recent sex (0/1)', '0.005', '(-0.003, 0.012)', '0.006', '(-0.003, 0.015)', '0.005', '(-0.002, 0.012)', '0.005', '(-0.003, 0.012)')
)

$cols = @('A', 'B', 'C', 'D', 'E', 'F', 'G', 'H', 'I')
$firstDataRow = 3

for ($r = 0; $r -lt $rows.Count; $r++) {
    $excelRow = $firstDataRow + $r
    $rowValues = $rows[$r]
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $cellRef = "$($cols[$c])$excelRow"
        $val = $rowValues[$c]
        $cell = $ws.Range($cellRef)

        # RD estimates like "0.011" / "-0.002" round-trip as plain numbers
        # unless the cell is pre-formatted as Text - the source table keeps
        # them as text (e.g. "-0.000", "0.048+"), so force that, then drop
        # back to the Normal style so no stray number-format residue stays
        # on the cell (values like "0.048+" already fail numeric parsing
        # and don't need this).
        $isNumericLooking = $val -match '^-?\d+\.\d+$'
        if ($isNumericLooking) {
            $cell.NumberFormat = "@"
            $cell.Value = $val
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }
}
